$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.65
$ws.Range("A6").Value = -22.291
$ws.Range("A7").Value = -21.286
$ws.Range("B7").Value = 6.312
$ws.Range("B12").Value = 4.853
$ws.Range("E12").Value = 17.268
$ws.Range("C13").Value = -13.225
$ws.Range("C14").Value = -12.152
$ws.Range("B15").Value = 5.087000000000001
$ws.Range("A16").Value = -21.949
$ws.Range("C16").Value = -13.14
$ws.Range("C19").Value = -12.228
$ws.Range("A20").Value = -21.433
$ws.Range("B20").Value = 6.121
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 9.415000000000001
$ws.Range("C22").Value = -12.689
$ws.Range("E22").Value = 16.974
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.767
$ws.Range("A29").Value = -21.267
$ws.Range("B29").Value = 6.122
$ws.Range("E29").Value = 17.075
$ws.Range("A32").Value = -21.671
$ws.Range("B34").Value = 7.603999999999999
$ws.Range("E34").Value = 16.638
$ws.Range("C36").Value = -12.37
$ws.Range("A40").Value = -20.452
$ws.Range("B42").Value = 7.582999999999998
$ws.Range("B43").Value = 5.795
$ws.Range("E43").Value = 17.192
$ws.Range("B44").Value = 5.096
$ws.Range("B45").Value = 5.048
$ws.Range("A46").Value = -21.15
$ws.Range("B46").Value = 6.344999999999999
$ws.Range("C46").Value = -13.898
$ws.Range("E48").Value = 17.205
$ws.Range("B50").Value = 5.291
$ws.Range("C50").Value = -13.363
$ws.Range("A51").Value = -20.724
$ws.Range("B51").Value = 7.513999999999998
$ws.Range("A52").Value = -21.262
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.467
$ws.Range("E60").Value = 16.552
$ws.Range("A62").Value = -21.782
$ws.Range("A66").Value = -21.213
$ws.Range("B66").Value = 5.829
$ws.Range("B67").Value = 5.249000000000001
$ws.Range("E68").Value = 17.148
$ws.Range("E70").Value = 17.596
$ws.Range("A73").Value = -20.129
$ws.Range("E73").Value = 16.544
$ws.Range("A74").Value = -21.045
$ws.Range("B79").Value = 5.423
$ws.Range("B84").Value = 5.798
$ws.Range("E87").Value = 16.403
$ws.Range("A92").Value = -21.239
$ws.Range("B92").Value = 5.548
$ws.Range("E92").Value = 18.022
$ws.Range("C95").Value = -11.644
$ws.Range("B97").Value = 6.49
$ws.Range("C97").Value = -12.765
$ws.Range("A100").Value = -21.69
$ws.Range("E101").Value = 16.6
